$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E values ("use leming, then stem" pass) for existing rows
$ws.Range("E2").Value = 64.739999999999995
$ws.Range("E3").Value = 64.400000000000006
$ws.Range("E4").Value = 67.040000000000006
$ws.Range("F4").Value = " (use leming, then stem)"
$ws.Range("E7").Value = 67.319999999999993
$ws.Range("E9").Value = 69.42

# Move selection to reflect the final active cell used while editing
$ws.Range("E9").Select()
